$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "yes"
$ws.Range("C22").Value = "yes"
$ws.Range("C26").Value = "yes"
$ws.Range("C27").Value = "yes"
$ws.Range("C30").Value = "yes"
$ws.Range("C37").Value = "yes"

$ws.Rows("17:21").AutoFit()

$ws.Range("C37").Select()
